# Auto-generated edit script applying the diff's cell-level changes.
$wb = $excel.ActiveWorkbook

# --- Sheet: 土地 ---
$ws = $wb.Worksheets.Item('土地')
$ws.Range("B1").Value = 'name'
$ws.Range("C1").Value = 'area'
$ws.Range("D1").Value = 'share_portion'
$ws.Range("E1").Value = 'owner'
$ws.Range("F1").Value = 'register_date'
$ws.Range("G1").Value = 'register_reason'
$ws.Range("H1").Value = 'acquire_value'
$ws.Range("I1").Value = 'property_category'
$ws.Range("J1").Value = 'category'
$ws.Range("K1").Value = 'date'
$ws.Range("L1").Value = 'legislator_name'
$ws.Range("M1").Value = 'legislator_id'
$ws.Range("N1").Value = 'source_file'
$ws.Range("O1").Value = 'index'
$ws.Range("B2").Value = '臺北市松山區敦化段三小段01690000地號'
$ws.Range("D2").Value = '325373分之5313'
$ws.Range("F2").Value = '89年10月04日'
$ws.Range("I2").Value = 'land'
$ws.Range("J2").Value = 'normal'
$ws.Range("K2").Value = '2013-12-24'
$ws.Range("L2").Value = '段宜康'
$ws.Range("M2").Value = 917
$ws.Range("N2").Value = 'tmpac2a1'
$ws.Range("O2").Value = 15
$ws.Range("B3").Value = '臺北市松山區敦化段三小段01710000地號'
$ws.Range("D3").Value = '325373分之5313'
$ws.Range("F3").Value = '89年10月04日'
$ws.Range("H3").Value = '''7095863'
$ws.Range("I3").Value = 'land'
$ws.Range("J3").Value = 'normal'
$ws.Range("K3").Value = '2013-12-24'
$ws.Range("L3").Value = '段宜康'
$ws.Range("M3").Value = 917
$ws.Range("N3").Value = 'tmpac2a1'
$ws.Range("O3").Value = 16
$ws.Range("B4").Value = '臺北市松山區敦化段三小段01710001地號'
$ws.Range("D4").Value = '325373分之5313'
$ws.Range("F4").Value = '89年10月04日'
$ws.Range("I4").Value = 'land'
$ws.Range("J4").Value = 'normal'
$ws.Range("K4").Value = '2013-12-24'
$ws.Range("L4").Value = '段宜康'
$ws.Range("M4").Value = 917
$ws.Range("N4").Value = 'tmpac2a1'
$ws.Range("O4").Value = 17
$ws.Range("B5").Value = '臺北市松山區敦化段三小段01910000地號'
$ws.Range("D5").Value = '325373分之5313'
$ws.Range("F5").Value = '89年10月04日'
$ws.Range("I5").Value = 'land'
$ws.Range("J5").Value = 'normal'
$ws.Range("K5").Value = '2013-12-24'
$ws.Range("L5").Value = '段宜康'
$ws.Range("M5").Value = 917
$ws.Range("N5").Value = 'tmpac2a1'
$ws.Range("O5").Value = 18
$ws.Range("B6").Value = '臺北市松山區延吉段三小段08320002地號'
$ws.Range("D6").Value = '10000分之302'
$ws.Range("F6").Value = '102年01月16曰'
$ws.Range("H6").Value = '4530635(無）'
$ws.Range("I6").Value = 'land'
$ws.Range("J6").Value = 'normal'
$ws.Range("K6").Value = '2013-12-24'
$ws.Range("L6").Value = '段宜康'
$ws.Range("M6").Value = 917
$ws.Range("N6").Value = 'tmpac2a1'
$ws.Range("O6").Value = 19

# --- Sheet: 建物 ---
$ws = $wb.Worksheets.Item('建物')
$ws.Range("B2").Value = '臺北市松山區敦化段三小段01482000建號'
$ws.Range("F2").Value = '89年10月04日'
$ws.Range("H2").Value = '''1269100'
$ws.Range("B3").Value = '臺北市松山區敦化段三小段02061000建號'
$ws.Range("D3").Value = '1000分之425'
$ws.Range("F3").Value = '89年10月04日'
$ws.Range("B4").Value = '臺北市松山區延吉段段三小段02423000建號'
$ws.Range("F4").Value = '102年01月16曰'
$ws.Range("H4").Value = '27426(無)'

# --- Sheet: 存款 ---
$ws = $wb.Worksheets.Item('存款')
$ws.Range("F2").Value = '''2220053'
$ws.Range("B3").Value = '合作金庫商業銀行光復南路分行'
$ws.Range("F3").Value = '''1659374'
$ws.Range("B5").Value = '台北富邦商業銀行敦南分行'
$ws.Range("B6").Value = '台北富邦商業銀行敦南分行'
$ws.Range("B7").Value = '台北富邦商業銀行敦南分行'
$ws.Range("F9").Value = '''1000342'

# --- Sheet: 保險 ---
$ws = $wb.Worksheets.Item('保險')
$ws.Range("B2").Value = '三商美邦人壽保險股份有限公司'
$ws.Range("B3").Value = '三商美邦人壽保險股份有限公司'
$ws.Range("D3").Value = '''30000'
$ws.Range("B8").Value = '一詮精密工業股份有限公司'
$ws.Range("D8").Value = '''30357'
$ws.Range("B9").Value = '太子建設開發股份有限公司'
$ws.Range("G9").Value = '''50000'

# --- Sheet: 具有相當價值之財產 ---
$ws = $wb.Worksheets.Item('具有相當價值之財產')
$ws.Range("C1").Value = '項件'
$ws.Range("B2").Value = '保險公司'
$ws.Range("C2").Value = '保險名稱'
$ws.Range("E2").Value = '備註'
